# Note de cadrage - "Ajout gestion emplacement dans les diagrammes"
#
# Current list (under "il faudra remplir la fiche article :") is:
#   Code produit
#   Hauteur, largeur, profondeur
#   Conditionnement
#       Lot 1 contenant 10 pieces
#   Reference du constructeur
#   Constructeur              <- carries the (hidden) "_GoBack" bookmark
#   Poids
#
# Target list:
#   Code produit
#   Hauteur, largeur, profondeur
#   Reference du constructeur
#   Constructeur
#   Poids
#   Emplacement                <- new bullet, now carries the "_GoBack" bookmark
#
# i.e. "Conditionnement" / "Lot 1 contenant 10 pieces" are removed, and a new
# "Emplacement" bullet is appended after "Poids", inheriting the bookmark that
# used to sit on "Constructeur".

$d = $word.ActiveDocument

function Find-ParagraphIndex($text) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.Trim() -eq $text) {
            return $i
        }
    }
    return -1
}

# --- 1. Drop "Conditionnement" and its sub-bullet "Lot 1 contenant 10 pieces" ---
$idxConditionnement = Find-ParagraphIndex "Conditionnement"
$idxLot = $idxConditionnement + 1

$pStart = $d.Paragraphs.Item($idxConditionnement)
$pEnd = $d.Paragraphs.Item($idxLot)
$toRemove = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$toRemove.Delete()

# --- 2. Move the hidden "_GoBack" bookmark off of "Constructeur" ---
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# --- 3. Insert the new "Emplacement" bullet right after "Poids" ---
$idxPoids = Find-ParagraphIndex "Poids"
$pPoids = $d.Paragraphs.Item($idxPoids)
$pPoids.Range.InsertParagraphAfter()

$idxEmplacement = $idxPoids + 1
$pEmplacement = $d.Paragraphs.Item($idxEmplacement)
# Append a throwaway trailing character so the bookmark range is non-empty
# (this interpreter only places a freshly-added bookmark correctly when the
# range it wraps contains visible text); deleting that character afterwards
# collapses the bookmark back onto the correct, now-empty, position instead
# of removing it - matching a bookmark that sits right after the run text.
$pEmplacement.Range.Text = "EmplacementX"

$rMarker = $d.Range($pEmplacement.Range.End - 2, $pEmplacement.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $rMarker)

$bmNew = $d.Bookmarks.Item("_GoBack")
$rExtra = $d.Range($bmNew.Start, $bmNew.End)
$rExtra.Delete()

Write-Output "done"
